$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 975.875
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 136
$ws.Range("H125").Value = 649.2857
$ws.Range("I125").Value = 515.8333
$ws.Range("K125").Value = 4642.4997
$ws.Range("M125").Value = -2182.4997
$ws.Range("H138").Value = 3145.4688
$ws.Range("I138").Value = 2633.9092
$ws.Range("K138").Value = 7901.7276
$ws.Range("M138").Value = -2761.7276
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3581.25
$ws.Range("I32").Value = 3666.8
$ws.Range("K32").Value = 3666.8
$ws.Range("M32").Value = -3379.8
$ws.Range("H44").Value = 25333.334
$ws.Range("I44").Value = 6000
$ws.Range("K44").Value = 6000
$ws.Range("M44").Value = -5512
$ws.Range("H74").Value = 4093.6843
$ws.Range("I74").Value = 3173.75
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 3173.75
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -2299.75
$ws.Range("N74").Value = -10748
$ws.Range("H77").Value = 4093.6843
$ws.Range("I77").Value = 3173.75
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 15868.75
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -11500.75
$ws.Range("N77").Value = -53736
$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 10000
$ws.Range("K137").Value = 10000
$ws.Range("M137").Value = -4900
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 806.4
$ws.Range("I86").Value = 833
$ws.Range("J86").Value = 700
$ws.Range("K86").Value = 833
$ws.Range("L86").Value = 700
$ws.Range("M86").Value = 290
$ws.Range("N86").Value = -2946
$ws.Range("H89").Value = 806.4
$ws.Range("I89").Value = 833
$ws.Range("J89").Value = 700
$ws.Range("K89").Value = 4165
$ws.Range("L89").Value = 3500
$ws.Range("M89").Value = 1451
$ws.Range("N89").Value = -14732
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 450
$ws.Range("J14").Value = 450
$ws.Range("L14").Value = 450
$ws.Range("N14").Value = -790
$ws.Range("H31").Value = 6398.4165
$ws.Range("J31").Value = 7159.857
$ws.Range("L31").Value = 7159.857
$ws.Range("N31").Value = -7749.857
$ws.Range("H34").Value = 6398.4165
$ws.Range("J34").Value = 7159.857
$ws.Range("L34").Value = 7159.857
$ws.Range("N34").Value = -7563.857
$ws.Range("H58").Value = 749.3333
$ws.Range("I58").Value = 749.3333
$ws.Range("K58").Value = 749.3333
$ws.Range("M58").Value = -546.3333
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H105").Value = 1783.4286
$ws.Range("I105").Value = 1196.8
$ws.Range("J105").Value = 3250
$ws.Range("K105").Value = 1196.8
$ws.Range("L105").Value = 3250
$ws.Range("M105").Value = 550.2
$ws.Range("N105").Value = -6744
$ws.Range("H136").Value = 749.3333
$ws.Range("I136").Value = 749.3333
$ws.Range("K136").Value = 2247.9999
$ws.Range("M136").Value = 302.0001000000002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 36.285713
$ws.Range("I6").Value = 54.75
$ws.Range("K6").Value = 164.25
$ws.Range("M6").Value = -51.25
$ws.Range("H10").Value = 2887.1428
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 4022
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 12066
$ws.Range("M10").Value = -11
$ws.Range("N10").Value = -12344
$ws.Range("H11").Value = 7143283.5
$ws.Range("I11").Value = 8333752
$ws.Range("K11").Value = 25001256
$ws.Range("M11").Value = -25001116
$ws.Range("H29").Value = 95
$ws.Range("I29").Value = 95
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 285
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -8
$ws.Range("N29").ClearContents()
$ws.Range("H55").Value = 1295.8823
$ws.Range("I55").Value = 961.8182
$ws.Range("J55").Value = 1908.3334
$ws.Range("K55").Value = 2885.4546
$ws.Range("L55").Value = 5725.0002
$ws.Range("M55").Value = -2708.4546
$ws.Range("N55").Value = -6079.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 428
$ws.Range("I9").Value = 537.3333
$ws.Range("K9").Value = 537.3333
$ws.Range("M9").Value = -367.3333
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970
$ws.Range("H102").Value = 1941.8572
$ws.Range("I102").Value = 1941.8572
$ws.Range("K102").Value = 1941.8572
$ws.Range("M102").Value = -319.8571999999999
$ws.Range("H126").Value = 3400
$ws.Range("I126").Value = 3400
$ws.Range("K126").Value = 10200
$ws.Range("M126").Value = -7730
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2038.8889
$ws.Range("I68").Value = 1981.25
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1981.25
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -1232.25
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 2038.8889
$ws.Range("I71").Value = 1981.25
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 9906.25
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -6162.25
$ws.Range("N71").Value = -19988
$ws.Range("H94").Value = 36731.8
$ws.Range("J94").Value = 38914.75
$ws.Range("L94").Value = 38914.75
$ws.Range("N94").Value = -40266.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 50000
$ws.Range("I114").Value = 50000
$ws.Range("K114").Value = 50000
$ws.Range("M114").Value = -45661
$ws.Range("H136").Value = 1849.8572
$ws.Range("I136").Value = 1941.5
$ws.Range("K136").Value = 5824.5
$ws.Range("M136").Value = -3274.5
